$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.481.16'
$ws.Range("E2").Value = '  +0.98%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.879.08'
$ws.Range("E3").Value = '  +0.82%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.08'
$ws.Range("E5").Value = '  +5.52%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4748'
$ws.Range("E7").Value = '  +1.59%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2897'
$ws.Range("E8").Value = '  +1.54%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06518'
$ws.Range("E9").Value = '  +0.55%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.62'
$ws.Range("E10").Value = '  +1.55%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07736'
$ws.Range("E11").Value = '  -0.40%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7430'
$ws.Range("E12").Value = '  +8.83%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '97.15'
$ws.Range("E13").Value = '  +3.47%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.878.66'
$ws.Range("E14").Value = '  -0.10%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.123'
$ws.Range("E15").Value = '  +1.45%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '273.79'
$ws.Range("E16").Value = '  +0.87%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.478.64'
$ws.Range("E17").Value = '  +1.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.65'
$ws.Range("E18").Value = '  +2.21%  '

# Row 19
$ws.Range("E19").Value = '  +0.00%  '

# Row 20
$ws.Range("E20").Value = '  +0.13%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.127.86'
$ws.Range("E21").Value = '  +0.39%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.10%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.265'
$ws.Range("E23").Value = '  +2.32%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.170'
$ws.Range("E24").Value = '  +0.86%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.278'
$ws.Range("E25").Value = '  -0.93%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.07'
$ws.Range("E26").Value = '  -0.86%  '

# Row 27
$ws.Range("E27").Value = '  +1.79%  '

# Row 28
$ws.Range("E28").Value = '  +3.10%  '

# Row 29
$ws.Range("E29").Value = '  +0.87%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09980'
$ws.Range("E30").Value = '  +1.45%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.519'
$ws.Range("E31").Value = '  +4.55%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.324'
$ws.Range("E32").Value = '  +1.96%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.058'
$ws.Range("E33").Value = '  +1.63%  '

# Row 34
$ws.Range("E34").Value = '  +2.24%  '

# Row 35
$ws.Range("E35").Value = '  +0.53%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6982'
$ws.Range("E36").Value = '  +1.26%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").Value = '  +0.40%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01868'
$ws.Range("E38").Value = '  +1.75%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.743'
$ws.Range("E39").Value = '  -0.34%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.355'
$ws.Range("E40").Value = '  -0.12%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '70.33'
$ws.Range("E41").Value = '  -0.95%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.922'
$ws.Range("E42").Value = '  +1.98%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4172'
$ws.Range("E43").Value = '  +2.71%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.11%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8364'
$ws.Range("E45").Value = '  +0.70%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.53'
$ws.Range("E46").Value = '  +0.11%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.271'
$ws.Range("E47").Value = '  +2.90%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.092'
$ws.Range("E48").Value = '  +1.80%  '

# Row 49
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '930.37'
$ws.Range("E49").Value = '  -0.44%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.40'
$ws.Range("E50").Value = '  +4.18%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05615'
$ws.Range("E51").Value = '  +0.69%  '
